# Fix typo: "Kraken2" -> "Kraken 2" throughout the "data" worksheet,
# and move the active cell selection from Z29 to Z3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Replace all occurrences of "Kraken2" with "Kraken 2" in the used range.
# This updates the shared string cell (Z2) as well as any literal text cells;
# cells containing formulas (Z3:Z25, which reference the previous row) will
# recalculate automatically and their cached values will follow.
$ws.Cells.Replace("Kraken2", "Kraken 2", 1, 1, $false, $false, $false, $false) | Out-Null

# Restore/update the active selection to Z3 (was Z29 before the edit).
$ws.Activate()
$ws.Range("Z3").Select() | Out-Null
